$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 ("lemmalist-greek") is removed entirely; rows below shift up.
$ws.Rows.Item(10).Delete()
